$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 627
$ws.Range("F3").Value = 629
$ws.Range("F4").Value = 911
$ws.Range("F5").Value = 678
$ws.Range("F6").Value = 807
$ws.Range("F7").Value = 374
$ws.Range("F8").Value = 578
$ws.Range("F10").Value = 1170
$ws.Range("F11").Value = 599
$ws.Range("F12").Value = 356
$ws.Range("F13").Value = 480
$ws.Range("F14").Value = 155
$ws.Range("F15").Value = 103
$ws.Range("F16").Value = 318
$ws.Range("F18").Value = 76
$ws.Range("F19").Value = 535
$ws.Range("F20").Value = 44
$ws.Range("F21").Value = 546
$ws.Range("F22").Value = 19
$ws.Range("F23").Value = 568

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 206
$ws.Range("F10").Value = 45

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 627
$ws.Range("F7").Value = 629
$ws.Range("F8").Value = 911
$ws.Range("F9").Value = 678
$ws.Range("F10").Value = 807
$ws.Range("F11").Value = 374
$ws.Range("F12").Value = 578
$ws.Range("F14").Value = 1170
$ws.Range("F15").Value = 599
$ws.Range("F18").Value = 356
$ws.Range("F19").Value = 480
$ws.Range("F21").Value = 155
$ws.Range("F22").Value = 103
$ws.Range("F24").Value = 318
$ws.Range("F26").Value = 76
$ws.Range("F27").Value = 206
$ws.Range("F28").Value = 45
$ws.Range("F29").Value = 535
$ws.Range("F33").Value = 44
$ws.Range("F34").Value = 546
$ws.Range("F35").Value = 19
$ws.Range("F36").Value = 568
